$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B11").NumberFormat = "General"
Write-Host "After: $($ws.Range('B11').NumberFormat)"
